# Add two new columns (I: "I0", J: "IF") to the sheet, matching the
# existing header style/format used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting from H1 (bold, centered, thin border) onto
# the two new header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-17 for columns I and J
$data = @{
    2  = @(5, 8)
    3  = @(4, 7)
    4  = @(6, 6)
    5  = @(7, 8)
    6  = @(8, 9)
    7  = @(7, 9)
    8  = @(7, 7)
    9  = @(7, 8)
    10 = @(9, 9)
    11 = @(8, 9)
    12 = @(5, 8)
    13 = @(1, 5)
    14 = @(6, 9)
    15 = @(6, 8)
    16 = @(1, 5)
    17 = @(6, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
